$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the 4 "section header" rows that referenced the old shared string
#    "LIVEHTA Automation - Test_NonOncology_Automation_1 - 3/6/2023" so that
#    they now reference the new text "LIVEHTA Automation -
#    Test_NonOncology_Automation_1". Once every usage of the old string is
#    gone it will naturally drop out of the shared string table.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "LIVEHTA Automation - Test_NonOncology_Automation_1"
$ws.Range("B4").Value = "LIVEHTA Automation - Test_NonOncology_Automation_1"
$ws.Range("B36").Value = "LIVEHTA Automation - Test_NonOncology_Automation_1"
$ws.Range("B57").Value = "LIVEHTA Automation - Test_NonOncology_Automation_1"

# ---------------------------------------------------------------------------
# 2. Column B got narrower (bestFit) because the text above got shorter.
#    Reproduce the new width as closely as this runtime's column-width
#    quantization allows.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 47.95

# ---------------------------------------------------------------------------
# 3. Append the new "pop5" block of rows (141-151).
# ---------------------------------------------------------------------------
$ws.Range("A141").Value = "pop5"
$ws.Range("B141").Value = "LIVEHTA Automation - Test_NonOncology_Automation_1"
$ws.Range("C141").Value = "File_with_ValidData.xlsx"
$ws.Range("D141").Value = "\Testdata\Non_Oncology\Templates\ImportPublications\File_with_ValidData.xlsx"
$ws.Range("E141").Value = 3
$ws.Range("E141").NumberFormat = "0"
$ws.Range("F141").Value = "Publication ID 1 has already been uploaded, please check and re-upload"

$ws.Range("A142").Value = "pop5"
$ws.Range("E142").Value = 11
$ws.Range("E142").NumberFormat = "0"
$ws.Range("F142").Value = "Publication ID 9 has already been uploaded, please check and re-upload"

$ws.Range("A143").Value = "pop5"
$ws.Range("E143").Value = 5
$ws.Range("E143").NumberFormat = "0"
$ws.Range("F143").Value = "Publication ID 3 has already been uploaded, please check and re-upload"

$ws.Range("A144").Value = "pop5"
$ws.Range("E144").Value = 13
$ws.Range("E144").NumberFormat = "0"
$ws.Range("F144").Value = "Publication ID 11 has already been uploaded, please check and re-upload"

$ws.Range("A145").Value = "pop5"
$ws.Range("E145").Value = 12
$ws.Range("E145").NumberFormat = "0"
$ws.Range("F145").Value = "Publication ID 10 has already been uploaded, please check and re-upload"

$ws.Range("A146").Value = "pop5"
$ws.Range("E146").Value = 4
$ws.Range("E146").NumberFormat = "0"
$ws.Range("F146").Value = "Publication ID 2 has already been uploaded, please check and re-upload"

$ws.Range("A147").Value = "pop5"
$ws.Range("E147").Value = 10
$ws.Range("E147").NumberFormat = "0"
$ws.Range("F147").Value = "Publication ID 8 has already been uploaded, please check and re-upload"

$ws.Range("A148").Value = "pop5"
$ws.Range("E148").Value = 9
$ws.Range("E148").NumberFormat = "0"
$ws.Range("F148").Value = "Publication ID 7 has already been uploaded, please check and re-upload"

$ws.Range("A149").Value = "pop5"
$ws.Range("E149").Value = 6
$ws.Range("E149").NumberFormat = "0"
$ws.Range("F149").Value = "Publication ID 4 has already been uploaded, please check and re-upload"

$ws.Range("A150").Value = "pop5"
$ws.Range("E150").Value = 8
$ws.Range("E150").NumberFormat = "0"
$ws.Range("F150").Value = "Publication ID 6 has already been uploaded, please check and re-upload"

$ws.Range("A151").Value = "pop5"
$ws.Range("E151").Value = 7
$ws.Range("E151").NumberFormat = "0"
$ws.Range("F151").Value = "Publication ID 5 has already been uploaded, please check and re-upload"

# ---------------------------------------------------------------------------
# 4. Update the view: the frozen pane now starts around row 127 and the
#    active selection moved from G1 to B146.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 127
$win.ScrollColumn = 1
$ws.Range("B146").Select()
